# 904. Fruit Into Baskets (Sliding window)
#
# Add a new LeetCode entry right after "219. Contains Duplicate II" in the
# "Sliding Window (Advance)" section. This inserts a row, which pushes every
# subsequent row (section headers, problems, hyperlink, summary formulas)
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 51 — same effect as right-clicking row 51 and choosing
# "Insert"; everything from the old row 51 onward shifts down by one and the
# new row inherits formatting from the row above it.
$ws.Rows.Item(51).Insert()

# Populate the new entry.
$ws.Range("A51").Value = 904
$ws.Range("B51").Value = "Fruit Into Baskets"
$ws.Range("C51").Value = "Medium"
$ws.Range("D51").Value = "Arrays,sliding window,hashmap"
$ws.Range("E51").Value = 45873
$ws.Range("F51").Value = "Python"

# The previously-selected cell shifts down from F51 to F52.
$ws.Range("F52").Select()

# The "127. Word Ladder" hyperlink lived on B129; its row moved to 130, so
# recreate the hyperlink there and restore the original (non-default) font
# formatting that Hyperlinks.Add() overwrites.
$ws.Range("B129").Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("B130"), "https://leetcode.com/problems/word-ladder/", "", "https://leetcode.com/problems/word-ladder/", "127. Word Ladder")
$ws.Range("B130").Font.Size = 10
$ws.Range("B130").Font.Color = 16745482
$ws.Range("B130").Font.Underline = $false
